# Update match-odds figures on the (single, active) sheet.
# These are plain numeric cell edits scattered across rows 4, 6, 8, 9 and 10 -
# no structural changes, no new rows/columns, just refreshed odds values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("H4").Value = 3.8
$ws.Range("M4").Value = 1.06
$ws.Range("N4").Value = 10
$ws.Range("Q4").Value = 2.05
$ws.Range("R4").Value = 1.75
$ws.Range("S4").Value = 1.4
$ws.Range("T4").Value = 2.75
$ws.Range("U4").Value = 2
$ws.Range("V4").Value = 1.73
$ws.Range("X4").Value = 7
$ws.Range("AB4").Value = 29
$ws.Range("AC4").Value = 9
$ws.Range("AE4").Value = 19
$ws.Range("AH4").Value = 15
$ws.Range("AJ4").Value = 19
$ws.Range("AS4").Value = 151
$ws.Range("AT4").Value = 2.75
$ws.Range("AU4").Value = 9
$ws.Range("AW4").Value = 7

# Row 6
$ws.Range("I6").Value = 4.5
$ws.Range("M6").Value = 1.06
$ws.Range("N6").Value = 10
$ws.Range("O6").Value = 1.33
$ws.Range("P6").Value = 3.25
$ws.Range("Q6").Value = 2.1
$ws.Range("R6").Value = 1.7
$ws.Range("X6").Value = 8
$ws.Range("AA6").Value = 15
$ws.Range("AB6").Value = 29
$ws.Range("AE6").Value = 17
$ws.Range("AN6").Value = 3.75
$ws.Range("AO6").Value = 9.5
$ws.Range("AS6").Value = 151
$ws.Range("AU6").Value = 8.5
$ws.Range("AZ6").Value = 81

# Row 8
$ws.Range("G8").Value = 2.35
$ws.Range("I8").Value = 2.85
$ws.Range("J8").Value = 2.95
$ws.Range("L8").Value = 3.45
$ws.Range("N8").Value = 7.8
$ws.Range("P8").Value = 3.6
$ws.Range("T8").Value = 2.72
$ws.Range("AA8").Value = 17.5
$ws.Range("AC8").Value = 7.8
$ws.Range("AJ8").Value = 10
$ws.Range("AK8").Value = 35
$ws.Range("AL8").Value = 22
$ws.Range("AM8").Value = 26
$ws.Range("AN8").Value = 4.35
$ws.Range("AT8").Value = 2.72
$ws.Range("AV8").Value = 55
$ws.Range("AW8").Value = 4.85
$ws.Range("AY8").Value = 22
$ws.Range("BA8").Value = 100
$ws.Range("BB8").Value = 250

# Row 9
$ws.Range("G9").Value = 1.95
$ws.Range("H9").Value = 3.55
$ws.Range("I9").Value = 3.55
$ws.Range("J9").Value = 2.52
$ws.Range("L9").Value = 3.95
$ws.Range("O9").Value = 1.23
$ws.Range("U9").Value = 1.62
$ws.Range("X9").Value = 10.25
$ws.Range("Y9").Value = 8.25
$ws.Range("Z9").Value = 17.5
$ws.Range("AA9").Value = 14.5
$ws.Range("AB9").Value = 22
$ws.Range("AD9").Value = 6.9
$ws.Range("AH9").Value = 12.5
$ws.Range("AI9").Value = 21
$ws.Range("AJ9").Value = 11.75
$ws.Range("AK9").Value = 50
$ws.Range("AL9").Value = 29
$ws.Range("AM9").Value = 32
$ws.Range("AN9").Value = 3.95
$ws.Range("AO9").Value = 10
$ws.Range("AP9").Value = 17
$ws.Range("AQ9").Value = 35
$ws.Range("AR9").Value = 60
$ws.Range("AU9").Value = 6.8
$ws.Range("AW9").Value = 5.5
$ws.Range("AX9").Value = 19
$ws.Range("AY9").Value = 24
$ws.Range("AZ9").Value = 100
$ws.Range("BA9").Value = 120
$ws.Range("BB9").Value = 300

# Row 10
$ws.Range("G10").Value = 2.6
$ws.Range("I10").Value = 2.62
$ws.Range("J10").Value = 3.15
$ws.Range("K10").Value = 2.07
$ws.Range("L10").Value = 3.25
$ws.Range("O10").Value = 1.3
$ws.Range("Q10").Value = 1.9
$ws.Range("X10").Value = 14
$ws.Range("Y10").Value = 9.5
$ws.Range("Z10").Value = 29
$ws.Range("AA10").Value = 20
$ws.Range("AB10").Value = 27
$ws.Range("AI10").Value = 13.5
$ws.Range("AJ10").Value = 9.75
$ws.Range("AK10").Value = 30
$ws.Range("AL10").Value = 22
$ws.Range("AM10").Value = 29
$ws.Range("AN10").Value = 4.55
$ws.Range("AO10").Value = 13.5
$ws.Range("AP10").Value = 20
$ws.Range("AQ10").Value = 60
$ws.Range("AW10").Value = 4.6
$ws.Range("AX10").Value = 14.5
$ws.Range("AY10").Value = 22
$ws.Range("AZ10").Value = 65
$ws.Range("BB10").Value = 300
